# Refresh "carjacking-by-neighborhood-by-month" report: roll the current
# "through" date from 2021-12-26 to 2021-12-27 (new data for 2022-01-04's
# pull), updating the sheet name / header label and the affected monthly
# counts (some rows get brand-new non-zero cells where the count was
# previously zero/blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title + the matching "current month" column header label.
$ws.Name = "Through 2021-12-27"
$ws.Range("B1").Value = "December 2021 (through December 27)"

# Updated / newly-populated monthly counts, by neighborhood row.
$ws.Range("BJ2").Value = 3     # West Town / December 2016
$ws.Range("N4").Value = 12     # North Lawndale / December 2020
$ws.Range("B5").Value = 3      # West Pullman / December 2021 (through Dec 27)
$ws.Range("AX6").Value = 9     # Garfield Park / December 2017
$ws.Range("B7").Value = 9      # Austin / December 2021 (through Dec 27)
$ws.Range("BV7").Value = 5     # Austin / December 2015
$ws.Range("AX8").Value = 2     # Chatham / December 2017
$ws.Range("B11").Value = 3     # Humboldt Park / December 2021 (through Dec 27)
$ws.Range("N12").Value = 3     # Little Italy, UIC / December 2020
$ws.Range("N13").Value = 6     # Roseland / December 2020
$ws.Range("N17").Value = 1     # United Center / December 2020 (new)
$ws.Range("AL17").Value = 1    # United Center / December 2018 (new)
$ws.Range("AL25").Value = 1    # Ashburn / December 2018 (new)
$ws.Range("AL26").Value = 3    # Lake View / December 2018
$ws.Range("B27").Value = 2     # Edgewater / December 2021 (through Dec 27)
$ws.Range("AX38").Value = 6    # Auburn Gresham / December 2017
$ws.Range("AL40").Value = 1    # Calumet Heights / December 2018 (new)
$ws.Range("B49").Value = 4     # Ukrainian Village / December 2021 (through Dec 27)
$ws.Range("AX66").Value = 2    # Avondale / December 2017
$ws.Range("AL82").Value = 2    # Logan Square / December 2018
$ws.Range("AX82").Value = 1    # Logan Square / December 2017 (new)
